$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211; this shifts the existing rows 211-284
# down to 212-285, matching the diff (dimension A1:R284 -> A1:R285).
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new record. The non-date /
# non-volume columns (K:R) repeat the values that were already present on
# the row directly below (now row 212), exactly as in the target sheet.
$ws.Range("A211").Value = 10
$ws.Range("B211").Value = "Vega Modelo de Temuco"
$ws.Range("C211").Value = "La Araucanía"
$ws.Range("D211").Value = 44559
$ws.Range("E211").Value = 9
$ws.Range("F211").Value = 100112040
$ws.Range("G211").Value = "Cilantro"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 110
$ws.Range("K211").Value = 6000
$ws.Range("L211").Value = 6000
$ws.Range("M211").Value = 6000
$ws.Range("N211").Value = "$/docena de atados (2 kilos)"
$ws.Range("O211").Value = "Provincia de Cautín"
$ws.Range("P211").Value = 3000
$ws.Range("Q211").Value = 2
$ws.Range("R211").Value = "Hortaliza"
